$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RunManager")

# Update the test data value for C4 from "Yes" to "No"
$ws.Range("C4").Value = "No"

# Update selection / active cell to C4 (as recorded in the saved view state)
$ws.Activate()
$ws.Range("C4").Select()
